# Daily attendance processing - 2025-10-11 12:32:41
# Reorders the comma-separated "Recorded By" names in column G for the
# affected session rows (no other data changes).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$updates = @{
    "G2"   = "system, System, backup@backdoor.com"
    "G3"   = "System, dnasr281@gmail.com"
    "G6"   = "System, dnasr281@gmail.com"
    "G7"   = "admin@admin.com, System"
    "G10"  = "System, dnasr281@gmail.com"
    "G11"  = "System, dnasr281@gmail.com"
    "G12"  = "System, dnasr281@gmail.com"
    "G13"  = "System, dnasr281@gmail.com"
    "G14"  = "System, dnasr281@gmail.com"
    "G15"  = "System, dnasr281@gmail.com"
    "G17"  = "System, dnasr281@gmail.com"
    "G29"  = "system, System, backup@backdoor.com"
    "G30"  = "System, dnasr281@gmail.com"
    "G33"  = "System, dnasr281@gmail.com"
    "G34"  = "admin@admin.com, System"
    "G37"  = "System, dnasr281@gmail.com"
    "G38"  = "System, dnasr281@gmail.com"
    "G39"  = "System, dnasr281@gmail.com"
    "G40"  = "System, dnasr281@gmail.com"
    "G41"  = "System, dnasr281@gmail.com"
    "G42"  = "System, dnasr281@gmail.com"
    "G44"  = "System, dnasr281@gmail.com"
    "G56"  = "system, System, backup@backdoor.com"
    "G57"  = "System, dnasr281@gmail.com"
    "G60"  = "System, dnasr281@gmail.com"
    "G61"  = "admin@admin.com, System"
    "G64"  = "System, dnasr281@gmail.com"
    "G65"  = "System, dnasr281@gmail.com"
    "G66"  = "System, dnasr281@gmail.com"
    "G67"  = "System, dnasr281@gmail.com"
    "G68"  = "System, dnasr281@gmail.com"
    "G69"  = "System, dnasr281@gmail.com"
    "G71"  = "System, dnasr281@gmail.com"
    "G86"  = "System, dnasr281@gmail.com"
    "G87"  = "System, dnasr281@gmail.com"
    "G88"  = "System, dnasr281@gmail.com"
    "G89"  = "System, dnasr281@gmail.com"
    "G90"  = "admin@admin.com, dnasr281@gmail.com"
    "G93"  = "System, dnasr281@gmail.com"
    "G95"  = "System, dnasr281@gmail.com"
    "G96"  = "System, dnasr281@gmail.com"
    "G112" = "System, dnasr281@gmail.com"
    "G113" = "System, dnasr281@gmail.com"
    "G114" = "System, dnasr281@gmail.com"
    "G115" = "System, dnasr281@gmail.com"
    "G116" = "admin@admin.com, dnasr281@gmail.com"
    "G119" = "System, dnasr281@gmail.com"
    "G121" = "System, dnasr281@gmail.com"
    "G122" = "System, dnasr281@gmail.com"
    "G138" = "System, dnasr281@gmail.com"
    "G139" = "System, dnasr281@gmail.com"
    "G140" = "System, dnasr281@gmail.com"
    "G141" = "System, dnasr281@gmail.com"
    "G142" = "admin@admin.com, dnasr281@gmail.com"
    "G145" = "System, dnasr281@gmail.com"
    "G147" = "System, dnasr281@gmail.com"
    "G148" = "System, dnasr281@gmail.com"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
